# Applies the diff: insert a new data row at row 40 (pushing existing
# rows 40-67 down to 41-68), and populate the new row 40 with the
# new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 40; this shifts rows 40..67
# down to 41..68 and extends the used range / dimension accordingly.
$ws.Rows.Item(40).Insert()

# Populate the newly inserted row 40 with the new record's data.
# Columns A,B,C,E,F,G,H,I,J,K,L,Q,T are constant across every data row
# in this sheet, so reuse those values for the new row too.
$ws.Range("A40").Value = 5
$ws.Range("B40").Value = "Macroferia Regional de Talca"
$ws.Range("C40").Value = "Maule"
$ws.Range("D40").Value = "2021-08-30"
$ws.Range("D40").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E40").Value = 7
$ws.Range("F40").Value = "Fruta"
$ws.Range("G40").Value = 100108
$ws.Range("H40").Value = "Tropicales y subtropicales"
$ws.Range("I40").Value = 100108002
$ws.Range("J40").Value = "Mango"
$ws.Range("K40").Value = "Sin especificar"
$ws.Range("L40").Value = "Primera"
$ws.Range("M40").Value = 102
$ws.Range("N40").Value = 8000
$ws.Range("O40").Value = 8000
$ws.Range("P40").Value = 8000
$ws.Range("Q40").Value = "$/bandeja 4 kilos"
$ws.Range("R40").Value = "Brasil"
$ws.Range("S40").Value = 2000
$ws.Range("T40").Value = 4
